$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55 (shifts existing rows 55.. down by one,
# including all formatting/styles carried from the surrounding rows).
$ws.Rows("55:55").Insert()

# Populate the newly inserted row 55 with the new weekly price-report entry.
$ws.Cells.Item(55, 1).Value = 11
$ws.Cells.Item(55, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value = "Bíobío"
$ws.Cells.Item(55, 4).Value = 45070
$ws.Cells.Item(55, 5).Value = 8
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100108
$ws.Cells.Item(55, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(55, 9).Value = 100108002
$ws.Cells.Item(55, 10).Value = "Mango"
$ws.Cells.Item(55, 11).Value = "Sin especificar"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 200
$ws.Cells.Item(55, 14).Value = 7500
$ws.Cells.Item(55, 15).Value = 8000
$ws.Cells.Item(55, 16).Value = 7750
$ws.Cells.Item(55, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(55, 18).Value = "Perú"
$ws.Cells.Item(55, 19).Value = 1938
$ws.Cells.Item(55, 20).Value = 4
